$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") for rows 2 through 46 moves from 45179 to 45180
# (one day later, stored as a date serial number).
for ($row = 2; $row -le 46; $row++) {
    $ws.Cells.Item($row, 3).Value = 45180
}
